$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 8510.683000000001
$ws.Range("I15").Value = 8510.683000000001
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 25532.049
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -25363.049

$ws.Range("H29").Value = 1000.6
$ws.Range("I29").Value = 1003
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 3009
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = -2728
$ws.Range("N29").Value = -3562

$ws.Range("H31").Value = 5933.3335
$ws.Range("I31").Value = 5950
$ws.Range("J31").Value = 5900
$ws.Range("K31").Value = 17850
$ws.Range("L31").Value = 17700
$ws.Range("M31").Value = -17620
$ws.Range("N31").Value = -18160

$ws.Range("H38").Value = 65.57143000000001
$ws.Range("I38").Value = 65.57143000000001
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 196.71429
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 175.28571

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H45").Value = 899
$ws.Range("I45").Value = 899
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2697
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2505

$ws.Range("H58").Value = 2088.3333
$ws.Range("I58").Value = 399.2857
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 1197.8571
$ws.Range("L58").Value = 24000
$ws.Range("M58").Value = -1047.8571

$ws.Range("H112").Value = 5478.657
$ws.Range("I112").Value = 2399.5
$ws.Range("J112").Value = 5665.273
$ws.Range("K112").Value = 7198.5
$ws.Range("L112").Value = 16995.819
$ws.Range("M112").Value = -6090.5
$ws.Range("N112").Value = -19211.819

$ws.Range("H116").Value = 886084.2
$ws.Range("I116").Value = 1355202.8
$ws.Range("J116").Value = 6486.75
$ws.Range("K116").Value = 1355202.8
$ws.Range("L116").Value = 6486.75
$ws.Range("M116").Value = -1351760.8
$ws.Range("N116").Value = -13370.75

$ws.Range("H132").Value = 35535.906
$ws.Range("I132").Value = 43004.766
$ws.Range("J132").Value = 3793.25
$ws.Range("K132").Value = 129014.298
$ws.Range("L132").Value = 11379.75
$ws.Range("M132").Value = -126484.298
$ws.Range("N132").Value = -16439.75

$ws.Range("H137").Value = 35030.79
$ws.Range("I137").Value = 28466.416
$ws.Range("J137").Value = 46284
$ws.Range("K137").Value = 85399.24800000001
$ws.Range("L137").Value = 138852
$ws.Range("M137").Value = -82849.24800000001
$ws.Range("N137").Value = -143952

$ws.Range("H138").Value = 68356.44
$ws.Range("I138").Value = 5370.4
$ws.Range("J138").Value = 173333.17
$ws.Range("K138").Value = 16111.2
$ws.Range("L138").Value = 519999.51
$ws.Range("M138").Value = -10971.2
$ws.Range("N138").Value = -530279.51

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14597.1875
$ws.Range("I32").Value = 14630.063
$ws.Range("J32").Value = 12000
$ws.Range("K32").Value = 14630.063
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = -14343.063

$ws.Range("H61").Value = 24452.777
$ws.Range("I61").Value = 12693.667
$ws.Range("J61").Value = 47971
$ws.Range("K61").Value = 12693.667
$ws.Range("L61").Value = 47971
$ws.Range("M61").Value = -12481.667
$ws.Range("N61").Value = -48395

$ws.Range("H74").Value = 140628.66
$ws.Range("I74").Value = 194768.61
$ws.Range("J74").Value = 11525.692
$ws.Range("K74").Value = 194768.61
$ws.Range("L74").Value = 11525.692
$ws.Range("M74").Value = -193894.61
$ws.Range("N74").Value = -13273.692

$ws.Range("H77").Value = 140628.66
$ws.Range("I77").Value = 194768.61
$ws.Range("J77").Value = 11525.692
$ws.Range("K77").Value = 973843.0499999999
$ws.Range("L77").Value = 57628.45999999999
$ws.Range("M77").Value = -969475.0499999999
$ws.Range("N77").Value = -66364.45999999999

$ws.Range("H136").Value = 24452.777
$ws.Range("I136").Value = 12693.667
$ws.Range("J136").Value = 47971
$ws.Range("K136").Value = 38081.001
$ws.Range("L136").Value = 143913
$ws.Range("M136").Value = -35531.001
$ws.Range("N136").Value = -149013

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 36571.43
$ws.Range("I141").Value = 36571.43
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 36571.43
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -31391.43
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1963315
$ws.Range("I31").Value = 3228243.2
$ws.Range("J31").Value = 2676.1
$ws.Range("K31").Value = 3228243.2
$ws.Range("L31").Value = 2676.1
$ws.Range("M31").Value = -3227948.2
$ws.Range("N31").Value = -3266.1

$ws.Range("H34").Value = 1963315
$ws.Range("I34").Value = 3228243.2
$ws.Range("J34").Value = 2676.1
$ws.Range("K34").Value = 3228243.2
$ws.Range("L34").Value = 2676.1
$ws.Range("M34").Value = -3228041.2
$ws.Range("N34").Value = -3080.1

$ws.Range("H132").Value = 63702.312
$ws.Range("I132").Value = 77586.84
$ws.Range("J132").Value = 3536
$ws.Range("K132").Value = 232760.52
$ws.Range("L132").Value = 10608
$ws.Range("M132").Value = -230230.52
$ws.Range("N132").Value = -15668

$ws.Range("H141").Value = 70000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 70000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3957.1836
$ws.Range("I68").Value = 1624.625
$ws.Range("J68").Value = 4412.317
$ws.Range("K68").Value = 4873.875
$ws.Range("L68").Value = 13236.951
$ws.Range("M68").Value = -4062.875
$ws.Range("N68").Value = -14858.951

$ws.Range("H71").Value = 3957.1836
$ws.Range("I71").Value = 1624.625
$ws.Range("J71").Value = 4412.317
$ws.Range("K71").Value = 14621.625
$ws.Range("L71").Value = 39710.853
$ws.Range("M71").Value = -10565.625
$ws.Range("N71").Value = -47822.853

$ws.Range("H107").Value = 4399.4
$ws.Range("I107").Value = 4399.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 13198.2
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -11278.2
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1416.717
$ws.Range("I126").Value = 1411.8914
$ws.Range("J126").Value = 1448.4286
$ws.Range("K126").Value = 4235.674199999999
$ws.Range("L126").Value = 4345.2858
$ws.Range("M126").Value = -1765.674199999999
$ws.Range("N126").Value = -9285.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4210.2383
$ws.Range("I46").Value = 1737.625
$ws.Range("J46").Value = 5731.846
$ws.Range("K46").Value = 1737.625
$ws.Range("L46").Value = 5731.846
$ws.Range("M46").Value = -1549.625
$ws.Range("N46").Value = -6107.846

$ws.Range("H55").Value = 2554.2
$ws.Range("I55").Value = 1747.5834
$ws.Range("J55").Value = 3764.125
$ws.Range("K55").Value = 1747.5834
$ws.Range("L55").Value = 3764.125
$ws.Range("M55").Value = -1574.5834
$ws.Range("N55").Value = -4110.125

$ws.Range("H132").Value = 6288.4
$ws.Range("I132").Value = 6111.125
$ws.Range("J132").Value = 6997.5
$ws.Range("K132").Value = 18333.375
$ws.Range("L132").Value = 20992.5
$ws.Range("M132").Value = -15803.375
$ws.Range("N132").Value = -26052.5

$ws.Range("H136").Value = 8121.5454
$ws.Range("I136").Value = 6792.25
$ws.Range("J136").Value = 11666.333
$ws.Range("K136").Value = 20376.75
$ws.Range("L136").Value = 34998.999
$ws.Range("M136").Value = -17826.75
$ws.Range("N136").Value = -40098.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 105000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 105000
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 105000
$ws.Range("N121").Value = -108494

$ws.Range("H126").Value = 174470.11
$ws.Range("I126").Value = 1776.4166
$ws.Range("J126").Value = 1003399.8
$ws.Range("K126").Value = 5329.2498
$ws.Range("L126").Value = 3010199.4
$ws.Range("M126").Value = -2859.2498
$ws.Range("N126").Value = -3015139.4

$ws.Range("H136").Value = 22108.844
$ws.Range("I136").Value = 25095.822
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 75287.466
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -72737.466
